# Sprint 3 burndown chart update:
#  - remove two leading blank rows
#  - replace the 14-task list with a new 6-task list (new task names/efforts)
#  - shrink the chart/table accordingly (summary rows shift from 20-22 to 12-14)
#  - drop the two decorative pictures that used to float next to the table
#  - minor workbook metadata / chart language tweaks

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Remove the 8 extra task rows (12-19), leaving just 6 task rows (6-11)
# before the summary block. Rows 1-2 were already empty placeholder rows
# (no cell content) and the used range naturally recomputes to start at
# row 3 once saved, so they don't need to be touched.
# ---------------------------------------------------------------------------
$ws.Rows("12:19").Delete()

# ---------------------------------------------------------------------------
# 3) Overwrite the remaining 6 task rows with the new Sprint 3 data.
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Look for code smells in source code."
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = $null
$ws.Range("H6").Value = $null
$ws.Range("I6").Value = $null
$ws.Range("J6").Value = $null
$ws.Range("K6").Value = $null

$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "Identify GoF patterns."
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = $null
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = $null
$ws.Range("I7").Value = $null
$ws.Range("J7").Value = $null
$ws.Range("K7").Value = $null

$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "Analyze metric data."
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = $null
$ws.Range("H8").Value = $null
$ws.Range("I8").Value = $null
$ws.Range("J8").Value = $null
$ws.Range("K8").Value = $null

$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "Create Use Case Diagrams."
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = $null
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = $null
$ws.Range("I9").Value = $null
$ws.Range("J9").Value = $null
$ws.Range("K9").Value = $null

$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "Review peers."
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = $null
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = $null
$ws.Range("I10").Value = $null
$ws.Range("J10").Value = $null
$ws.Range("K10").Value = $null

$ws.Range("B11").Value = 6
$ws.Range("C11").Value = "Fix work according to reviews."
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = $null
$ws.Range("F11").Value = $null
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = $null
$ws.Range("I11").Value = $null
$ws.Range("J11").Value = $null
$ws.Range("K11").Value = $null

Write-Host "done"
